# Daily attendance processing - 2026-01-26 14:43:57
# Reorder the semicolon/comma separated "Recorded By" lists in column G
# to reflect the latest recorded-by ordering for each session row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    "G2"  = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg"
    "G3"  = "Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    "G4"  = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    "G5"  = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    "G6"  = "alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    "G7"  = "lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
    "G11" = "aya.saeed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
    "G12" = "Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
    "G13" = "yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
    "G15" = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    "G17" = "esraa.sami@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    "G19" = "Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
    "G20" = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    "G24" = "Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
    "G27" = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
    "G30" = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
